# Add a new "2023" data column (Q) to the existing Kyrgyz SDG indicator
# table, mirroring the formatting of the last existing year column (P),
# and adjust the row heights that Excel re-flowed when the new column
# was added.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Clone formatting from column P (2022) into the new column Q ---
# so the new cells pick up the same styles (borders, number format,
# fonts, alignment) as the rest of the table, without creating any new
# cell-style entries.
$ws.Range("P3:P14").Copy()
$ws.Range("Q3:Q14").PasteSpecial(-4122)

# --- 2. Populate the new 2023 values ---
$ws.Range("Q4").Value = 2023
$ws.Range("Q5").Value = 74.605426356589135
$ws.Range("Q6").Value = 118.8
$ws.Range("Q7").Value = 71.61643835616438
$ws.Range("Q8").Value = 95.703125
$ws.Range("Q9").Value = 113.91018619934282
$ws.Range("Q10").Value = 108.21501014198785
$ws.Range("Q11").Value = 165.26684164479443
$ws.Range("Q12").Value = 48.504446240905416
$ws.Range("Q13").Value = 97.361348644026393
$ws.Range("Q14").Value = 52.747252747252752

# Q3 stays empty (just the bottom border style, like N3:P3).

# --- 3. Row heights changed after the new column/data was laid out ---
$ws.Rows.Item(4).RowHeight = 16.5
$ws.Rows.Item(5).RowHeight = 27
$ws.Rows.Item(6).RowHeight = 24.75
$ws.Rows.Item(7).RowHeight = 16.5
$ws.Rows.Item(8).RowHeight = 16.5
$ws.Rows.Item(9).RowHeight = 16.5
$ws.Rows.Item(10).RowHeight = 16.5
$ws.Rows.Item(11).RowHeight = 16.5
$ws.Rows.Item(12).RowHeight = 16.5
$ws.Rows.Item(13).RowHeight = 16.5
$ws.Rows.Item(14).RowHeight = 16.5

# --- 4. Reset the selection back to the top-left cell ---
# (the saved file originally had a stray selection at R1; put the
# cursor back on A1 now that the used range extends to column Q)
$ws.Range("A1").Select()
